$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume Number + report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/14/2025  Through  4/20/2025"

# --- Weekly crime statistics table updates (rows 14-28) ---
# Row 14
$ws.Range("N14").Value = -83.333333333333

# Row 15
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -76.923076923076

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = 0
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = 73.333333333333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -72.916666666666
$ws.Range("N16").Value = -91.304347826087

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 65
$ws.Range("J17").Value = 86
$ws.Range("K17").Value = -24.418604651162
$ws.Range("L17").Value = 1.5625
$ws.Range("M17").Value = -26.136363636363
$ws.Range("N17").Value = -31.578947368421

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 63.636363636363
$ws.Range("I18").Value = 51
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = 30.769230769230
$ws.Range("L18").Value = 13.333333333333
$ws.Range("M18").Value = -52.336448598130
$ws.Range("N18").Value = -89.079229122055

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 66.666666666666
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 21.428571428571
$ws.Range("I19").Value = 114
$ws.Range("J19").Value = 124
$ws.Range("K19").Value = -8.064516129032
$ws.Range("L19").Value = 3.636363636363
$ws.Range("M19").Value = -20.279720279720
$ws.Range("N19").Value = -27.848101265822

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -90
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = -41.666666666666
$ws.Range("I20").Value = 54
$ws.Range("J20").Value = 69
$ws.Range("K20").Value = -21.739130434782
$ws.Range("L20").Value = 1.886792452830
$ws.Range("M20").Value = -51.785714285714
$ws.Range("N20").Value = -94.871794871794

# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -16.666666666666
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = 7.058823529411
$ws.Range("I21").Value = 314
$ws.Range("J21").Value = 339
$ws.Range("K21").Value = -7.374631268436
$ws.Range("L21").Value = 4.318936877076
$ws.Range("M21").Value = -43.727598566308
$ws.Range("N21").Value = -84.983261597321

# Row 24
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 30
$ws.Range("G24").Value = 43
$ws.Range("H24").Value = -30.232558139534
$ws.Range("I24").Value = 158
$ws.Range("J24").Value = 207
$ws.Range("K24").Value = -23.671497584541
$ws.Range("L24").Value = -32.765957446808
$ws.Range("M24").Value = -33.333333333333

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 12.5
$ws.Range("I25").Value = 37
$ws.Range("J25").Value = 47
$ws.Range("K25").Value = -21.276595744680
$ws.Range("L25").Value = -9.756097560975

# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 47
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 161.111111111111
$ws.Range("I26").Value = 134
$ws.Range("J26").Value = 95
$ws.Range("K26").Value = 41.052631578947
$ws.Range("L26").Value = 42.553191489361
$ws.Range("M26").Value = -29.842931937172

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 13
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = -13.333333333333
$ws.Range("L28").Value = 225
